$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (theta_se) standard-error values, columns B..J
$ws.Range("B4").Value = "(0.42)"
$ws.Range("C4").Value = "(0.49)"
$ws.Range("D4").Value = "(0.07)"
$ws.Range("E4").Value = "(0.72)"
$ws.Range("F4").Value = "(0.48)"
$ws.Range("G4").Value = "(0.12)"
$ws.Range("H4").Value = "(0.13)"
$ws.Range("I4").Value = "(0.57)"
$ws.Range("J4").Value = "(0.87)"

# Row 6 (lambda_se) standard-error values, columns B..J
$ws.Range("B6").Value = "(0.39)"
$ws.Range("C6").Value = "(0.39)"
$ws.Range("D6").Value = "(0.26)"
$ws.Range("E6").Value = "(0.51)"
$ws.Range("F6").Value = "(0.18)"
$ws.Range("G6").Value = "(0.05)"
$ws.Range("H6").Value = "(0.16)"
$ws.Range("I6").Value = "(0.29)"
$ws.Range("J6").Value = "(0.52)"
